$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 700, shifting existing rows 700-765 down to 701-766
$ws.Rows.Item(700).Insert()

# Populate the new row 700 with the new record's data
$ws.Range("A700").Value = 3
$ws.Range("B700").Value = "Femacal de La Calera"
$ws.Range("C700").Value = "Coquimbo"
$ws.Range("D700").Value = 45166
$ws.Range("E700").Value = 5
$ws.Range("F700").Value = 100112021
$ws.Range("G700").Value = "Ají"
$ws.Range("H700").Value = "Inferno"
$ws.Range("I700").Value = "Primera"
$ws.Range("J700").Value = 68
$ws.Range("K700").Value = 15000
$ws.Range("L700").Value = 16000
$ws.Range("M700").Value = 15441
$ws.Range("N700").Value = "$/caja 10 kilos"
$ws.Range("O700").Value = "Región de Arica y Parinacota"
$ws.Range("P700").Value = 1544
$ws.Range("Q700").Value = 10
$ws.Range("R700").Value = "Hortaliza"
